$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column (J) values for the practice rows (pair_kind = "generic")
$ws.Range("J2").Value2 = "generic"
$ws.Range("J3").Value2 = "generic"
$ws.Range("J4").Value2 = "generic"
$ws.Range("J5").Value2 = "generic"

# New block of "stim details" data added near the bottom of the sheet
$ws.Range("A27").Value2 = "stim details"

$ws.Range("A28").Value2 = "month"
$ws.Range("B28").Value2 = "word_type"
$ws.Range("C28").Value2 = "need_audio"
$ws.Range("D28").Value2 = "need_image"
$ws.Range("E28").Value2 = "word"
$ws.Range("F28").Value2 = "count"
$ws.Range("G28").Value2 = "find images"

$ws.Range("A29").Value2 = 6
$ws.Range("B29").Value2 = "video"

$ws.Range("A30").Value2 = 6
$ws.Range("B30").Value2 = "video"

$ws.Range("A31").Value2 = 7
$ws.Range("B31").Value2 = "video"

$ws.Range("A32").Value2 = 7
$ws.Range("B32").Value2 = "video"

$ws.Range("A33").Value2 = 6
$ws.Range("B33").Value2 = "audio"

$ws.Range("A34").Value2 = 6
$ws.Range("B34").Value2 = "audio"

$ws.Range("A35").Value2 = 7
$ws.Range("B35").Value2 = "audio"

$ws.Range("A36").Value2 = 7
$ws.Range("B36").Value2 = "audio"
